$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.250.24'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '3.672.79'
$ws.Range('E3').Value = '  -0.40%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '674.97'
$ws.Range('E5').Value = '  -0.95%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.11'
$ws.Range('E6').Value = '  -2.74%  '
$ws.Range('E7').Value = '  -0.09%  '
$ws.Range('E8').Value = '  -1.38%  '
$ws.Range('E9').Value = '  -1.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.95'
$ws.Range('E10').Value = '  -5.39%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.435'
$ws.Range('E11').Value = '  -2.76%  '
$ws.Range('E12').Value = '  -3.42%  '
$ws.Range('D13').Value = '4.293.73'
$ws.Range('E13').Value = '  -0.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.29'
$ws.Range('E14').Value = '  -3.86%  '
$ws.Range('D15').Value = '3.669.55'
$ws.Range('E15').Value = '  -0.41%  '
$ws.Range('D16').Value = '69.190.81'
$ws.Range('E16').Value = '  -0.33%  '
$ws.Range('E17').Value = '  +1.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '16.04'
$ws.Range('E18').Value = '  -1.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.42'
$ws.Range('E19').Value = '  -3.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '466.28'
$ws.Range('E20').Value = '  -3.29%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.97'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.648'
$ws.Range('E22').Value = '  -2.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '79.71'
$ws.Range('E23').Value = '  -0.73%  '
$ws.Range('D24').Value = '3.819.04'
$ws.Range('E24').Value = '  -0.39%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('E26').Value = '  -6.70%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.91'
$ws.Range('E27').Value = '  -5.08%  '
$ws.Range('E28').Value = '  -5.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.67'
$ws.Range('E29').Value = '  -1.31%  '
$ws.Range('E30').Value = '  -5.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.60'
$ws.Range('E31').Value = '  -3.99%  '
$ws.Range('E32').Value = '  +0.59%  '
$ws.Range('E33').Value = '  -5.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.84'
$ws.Range('E34').Value = '  -1.09%  '
$ws.Range('D35').Value = '3.666.53'
$ws.Range('E35').Value = '  +0.22%  '
$ws.Range('E36').Value = '  -5.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.17'
$ws.Range('E37').Value = '  -3.77%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.23'
$ws.Range('E38').Value = '  -1.55%  '
$ws.Range('E40').Value = '  -0.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.21'
$ws.Range('E41').Value = '  -2.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '174.13'
$ws.Range('E42').Value = '  +8.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0898'
$ws.Range('E43').Value = '  -4.28%  '
$ws.Range('E44').Value = '  -1.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '47.53'
$ws.Range('E45').Value = '  -1.87%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '28.08'
$ws.Range('E46').Value = '  -6.63%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.69'
$ws.Range('E47').Value = '  -5.11%  '
$ws.Range('E48').Value = '  -4.55%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.29'
$ws.Range('E49').Value = '  -4.37%  '
$ws.Range('E50').Value = '  -4.38%  '
$ws.Range('E51').Value = '  -3.14%  '
